$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3
$ws.Range("B3").Value = 551695.774015071

# Clear D3 (becomes empty / inline string cell)
$ws.Range("D3").Value = ""

# Update C4
$ws.Range("C4").Value = 11.92123037940973

# Update C5
$ws.Range("C5").Value = 0

# Row 7: rename "Other" -> "Biogas", update D7
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 45.06162629812869

# New Row 8: "Other" with D8 value, copying style from row 7
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 170.1874117528692
